$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Lvl0")
$ws2 = $wb.Worksheets.Item("Lvl1")
$ws3 = $wb.Worksheets.Item("Lvl2")

# ---------------------------------------------------------------
# Sheet1 (Lvl0): move the trailing "test2/TypeZ" example block down
# (rows 18-19 -> rows 25-26) to make room for the two new examples
# (seasonEmptyArr / seasonArrWithEmptyItems) inserted at rows 17-18.
# ---------------------------------------------------------------
$ws1.Range("B25").Value = $ws1.Range("B18").Value2

$ws1.Range("C26").Value = $ws1.Range("C19").Value2
$ws1.Range("D26").Value = $ws1.Range("D19").Value2
$ws1.Range("E26").Value = $ws1.Range("E19").Value2
$ws1.Range("F26").Value = $ws1.Range("F19").Value2

$ws1.Range("B18").ClearContents()
$ws1.Range("C19").ClearContents()
$ws1.Range("D19").ClearContents()
$ws1.Range("E19").ClearContents()
$ws1.Range("F19").ClearContents()

# New row 17: field5 / arr / seasonEmptyArr (an array field with zero items)
$ws1.Range("D17").Value = "field5"
$ws1.Range("E17").Value = "arr"
$ws1.Range("F17").Value = "seasonEmptyArr"

# New row 18: field6 / arr / seasonArrWithEmptyItems (an array field whose
# items carry no data of their own beyond their id)
$ws1.Range("D18").Value = "field6"
$ws1.Range("E18").Value = "arr"
$ws1.Range("F18").Value = "seasonArrWithEmptyItems"

$ws1.Columns.Item(6).ColumnWidth = 27.6

# ---------------------------------------------------------------
# Sheet2 (Lvl1): add the definitions referenced from Lvl0.
# ---------------------------------------------------------------
$ws2.Range("C17").Value = "seasonEmptyArr"

$ws2.Range("C20").Value = "seasonArrWithEmptyItems"
$ws2.Range("E20").Value = "ref "
$ws2.Range("F20").Value = "seasonArrWithEmptyItems_item1"

$ws2.Range("E21").Value = "ref "
$ws2.Range("F21").Value = "seasonArrWithEmptyItems_item2"

$ws2.Columns.Item(3).ColumnWidth = 27.6

# ---------------------------------------------------------------
# Sheet3 (Lvl2): the two (empty) array items themselves - only the id.
# ---------------------------------------------------------------
$ws3.Range("C15").Value = "seasonArrWithEmptyItems_item1"
$ws3.Range("C16").Value = "seasonArrWithEmptyItems_item2"

$ws3.Columns.Item(3).ColumnWidth = 36.6

# ---------------------------------------------------------------
# Selections / active sheet - Lvl0 becomes the active tab, with the
# cursor left on the newly added "seasonArrWithEmptyItems" cell.
# ---------------------------------------------------------------
$ws2.Activate()
$ws2.Range("F20:F21").Select()

$ws3.Activate()
$ws3.Range("G9").Select()

$ws1.Activate()
$ws1.Range("F18").Select()
